$wb = $excel.ActiveWorkbook

# Insert a brand new worksheet as the very first tab, which will hold the
# 2021 pollution data. Worksheets.Add(Before) puts the new sheet directly
# before the sheet passed in, so anchoring on the current first tab puts
# our new sheet at position 1 and shifts everything else down by one.
$firstSheet = $wb.Worksheets.Item(1)
$ws2021 = $wb.Worksheets.Add($firstSheet)
$ws2021.Name = "data2021"

# Header row.
$ws2021.Cells.Item(1,1).Value = "woj."
$ws2021.Cells.Item(1,2).Value = "area"
$ws2021.Cells.Item(1,3).Value = "co2"
$ws2021.Cells.Item(1,4).Value = "metan"
$ws2021.Cells.Item(1,5).Value = "n2o"
$ws2021.Cells.Item(1,6).Value = "so2"
$ws2021.Cells.Item(1,7).Value = "no"
$ws2021.Cells.Item(1,8).Value = "co"

# Data rows (2021 values, one Polish voivodeship per row).
$data = @(
    @("dolnoslaskie",         19947, 24306.62,              42.2,                 3.18,  26.56, 37.65,  158.11000000000001),
    @("kujawsko-pomorskie",   17972, 15750.88,              62.3,                 6.2,   27.18, 44.77,  175.51),
    @("lubelskie",            25122, 12270.73,              130.4,                6.11,  19.91, 31.41,  145.78),
    @("lubuskie",             13988, 4986.79,                40.200000000000003,  1.68,  3.98,  12.54,  61.75),
    @("lodzkie",              18219, 48236.01,               93.93,               5.26,  67.25, 66.62,  229.63),
    @("malopolskie",          15183, 17924.34,               45.04,               2.35,  22.78, 32.44,  156.94999999999999),
    @("mazowieckie",          35558, 58045.52,              148.77000000000001,   9.9,   49.15, 87.88,  229.28),
    @("opolskie",              9412, 20797.11,                23.12,              2.76,  13.43, 25.71,   68.319999999999993),
    @("podkarpackie",         17846, 8028.31,                 31.42,              1.54,  16.39, 19.260000000000002, 242.66),
    @("podlaskie",            20187, 4119.51,                  0,                 7.32,   7.16, 17.600000000000001,  68.64),
    @("pomorskie",            18310, 13564.65,                53.3,               3.92,  23.21, 35,     161.24),
    @("slaskie",              12333, 50417.5,                488.88,              2.61,  44.31, 55.84,  219.66),
    @("swietokrzyskie",       11711, 15333.41,                22.97,              1.87,  16.02, 22.2,    68.98),
    @("warminsko-mazurskie",  24173, 5524.27,                 58.35,              5.67,  12.09, 21.25,    0),
    @("wielkopolskie",        29826, 21070.68,               141.33000000000001, 11.73,  23.5,  50.53,  213.79),
    @("zachodniopomorskie",   22892, 11199.94,                38.89,              5.69,  19.47, 30.68,  195.65)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($value in $row) {
        $ws2021.Cells.Item($r, $c).Value = $value
        $c = $c + 1
    }
    $r = $r + 1
}

# Match the author's final selection on the new sheet.
$null = $ws2021.Range("F21").Select()
